$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New bug report row (row 2)
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "View dispatch doesn't work"
$ws.Range("C2").Value = "OPEN"
$ws.Range("D2").Value = "Jobs -> Dispatches"
$ws.Range("E2").Value = 'The document created when pressing the "View dispatch" button is just a static sample document rather than content generated from the actual dispatch.'
$ws.Range("F2").Value = 40245
$ws.Range("G2").Value = 40245

# Match the row height used for the wrapped description text
$ws.Rows.Item(2).RowHeight = 60

# Move the active selection the way the author left it
[void]$ws.Range("F3").Select()
